$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J3").Value = "'10/8/2022"
$ws.Range("J3").NumberFormat = "mm-dd-yy"

$ws.Range("J6").Value = "'10/15/2022"
$ws.Range("J6").NumberFormat = "mm-dd-yy"

$ws.Range("J8").Value = "'10/29/2022"
$ws.Range("J8").NumberFormat = "mm-dd-yy"

$ws.Range("J11").Value = "'10/4/2022"
$ws.Range("J11").NumberFormat = "mm-dd-yy"

$ws.Range("J13").Value = "'11/12/2022"
$ws.Range("J13").NumberFormat = "mm-dd-yy"

$ws.Range("J15").Value = "'11/26/2022"
$ws.Range("J15").NumberFormat = "mm-dd-yy"
